$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 20 (shifts old row 20 "删除用例/deletecass" down to row 21,
# and expands the "用例管理" merge A16:A20 to A16:A21 automatically).
$ws.Rows.Item(20).Insert()

# Fill in the new "执行用例" (runtests) interface row.
$ws.Range("C20").Value = "runtests"
$ws.Range("B20").Value = "执行用例"
$ws.Range("D20").Value = "post"
$ws.Range("E20").Value = "{""idlist"":""1,2""}"
$ws.Range("F19").Copy($ws.Range("F20"))

# Reflect the final view/selection state (best effort: engine has no
# window-scroll object, so we at least move the active cell to F20).
$ws.Range("F20").Select() | Out-Null
